$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1965").Value = 'Buying Opportunity'
$ws.Range("B1965").Value = 'support Zone'
$ws.Range("C1965").Value = 'long buildup'
$ws.Range("D1965").Value = 'Short buildup'
$ws.Range("E1965").Value = 'FII ENTERING'

$ws.Range("A1966").Value = '3MINDIA'
$ws.Range("B1966").Value = 'BSE'
$ws.Range("F1966").Value = 39600.65
$ws.Range("G1966").Value = 2412.3

$ws.Range("A1967").Value = 'AARTIIND'
$ws.Range("B1967").Value = 'CYBERTECH'
$ws.Range("F1967").Value = 742.8
$ws.Range("G1967").Value = 195.35

$ws.Range("A1968").Value = 'ACEINTEG'
$ws.Range("B1968").Value = 'DONEAR'
$ws.Range("F1968").Value = 36.97
$ws.Range("G1968").Value = 138

$ws.Range("A1969").Value = 'ADSL'
$ws.Range("B1969").Value = 'HTMEDIA'
$ws.Range("F1969").Value = 233.26
$ws.Range("G1969").Value = 26.59

$ws.Range("A1970").Value = 'AJANTPHARM'
$ws.Range("B1970").Value = 'INVENTURE'
$ws.Range("F1970").Value = 2530.55
$ws.Range("G1970").Value = 2.9

$ws.Range("A1971").Value = 'ANDHRAPAP'
$ws.Range("B1971").Value = 'KOTHARIPET'
$ws.Range("F1971").Value = 582.1
$ws.Range("G1971").Value = 173.78

$ws.Range("A1972").Value = 'ANDHRSUGAR'
$ws.Range("B1972").Value = 'LOTUSEYE'
$ws.Range("F1972").Value = 120.78
$ws.Range("G1972").Value = 66.94

$ws.Range("A1973").Value = 'APARINDS'
$ws.Range("B1973").Value = 'MHRIL'
$ws.Range("F1973").Value = 9109.950000000001
$ws.Range("G1973").Value = 457

$ws.Range("A1974").Value = 'ASTERDM'
$ws.Range("B1974").Value = 'MURUDCERA'
$ws.Range("F1974").Value = 348.05
$ws.Range("G1974").Value = 57.12

$ws.Range("A1975").Value = 'BASML'
$ws.Range("B1975").Value = 'OSWALGREEN'
$ws.Range("F1975").Value = 59.12
$ws.Range("G1975").Value = 36.93

$ws.Range("A1976").Value = 'CHEVIOT'
$ws.Range("B1976").Value = 'PAR'
$ws.Range("F1976").Value = 1412.6
$ws.Range("G1976").Value = 242.54

$ws.Range("A1977").Value = 'COCHINSHIP'
$ws.Range("B1977").Value = 'PNC'
$ws.Range("F1977").Value = 2678.25
$ws.Range("G1977").Value = 69.95

$ws.Range("A1978").Value = 'COLPAL'
$ws.Range("B1978").Value = 'RAYMOND'
$ws.Range("F1978").Value = 3366.55
$ws.Range("G1978").Value = 2035.35

$ws.Range("A1979").Value = 'DIXON'
$ws.Range("F1979").Value = 11977.35

$ws.Range("A1980").Value = 'EDELWEISS'
$ws.Range("F1980").Value = 68.56

$ws.Range("A1981").Value = 'FCL'
$ws.Range("F1981").Value = 385.25

$ws.Range("A1982").Value = 'GLAND'
$ws.Range("F1982").Value = 2058.65

$ws.Range("A1983").Value = 'GOKULAGRO'
$ws.Range("F1983").Value = 183.21

$ws.Range("A1984").Value = 'GREENPANEL'
$ws.Range("F1984").Value = 356.2

$ws.Range("A1985").Value = 'HARRMALAYA'
$ws.Range("F1985").Value = 258.55

$ws.Range("A1986").Value = 'HDFCSML250'
$ws.Range("F1986").Value = 180.65

$ws.Range("A1987").Value = 'HFCL'
$ws.Range("F1987").Value = 134.81

$ws.Range("A1988").Value = 'HINDPETRO'
$ws.Range("F1988").Value = 395.75

$ws.Range("A1989").Value = 'HIRECT'
$ws.Range("F1989").Value = 762.45

$ws.Range("A1990").Value = 'IGARASHI'
$ws.Range("F1990").Value = 595.85

$ws.Range("A1991").Value = 'IIFLSEC'
$ws.Range("F1991").Value = 214.2

$ws.Range("A1992").Value = 'INGERRAND'
$ws.Range("F1992").Value = 4468.6

$ws.Range("A1993").Value = 'INOXGREEN'
$ws.Range("F1993").Value = 182.42

$ws.Range("A1994").Value = 'IPL'
$ws.Range("F1994").Value = 225.54

$ws.Range("A1995").Value = 'ISGEC'
$ws.Range("F1995").Value = 1548.25

$ws.Range("A1996").Value = 'ITDCEM'
$ws.Range("F1996").Value = 524.45

$ws.Range("A1997").Value = 'JINDALSAW'
$ws.Range("F1997").Value = 635.65

$ws.Range("A1998").Value = 'JKLAKSHMI'
$ws.Range("F1998").Value = 900.35

$ws.Range("A1999").Value = 'KABRAEXTRU'
$ws.Range("F1999").Value = 426.65

$ws.Range("A2000").Value = 'KANSAINER'
$ws.Range("F2000").Value = 304.55

$ws.Range("A2001").Value = 'LAOPALA'
$ws.Range("F2001").Value = 323.65

$ws.Range("A2002").Value = 'MANOMAY'
$ws.Range("F2002").Value = 280.64

$ws.Range("A2003").Value = 'MANYAVAR'
$ws.Range("F2003").Value = 1120.45

$ws.Range("A2004").Value = 'MEGASOFT'
$ws.Range("F2004").Value = 73.69

$ws.Range("A2005").Value = 'NAHARPOLY'
$ws.Range("F2005").Value = 325.36

$ws.Range("A2006").Value = 'NAVINFLUOR'
$ws.Range("F2006").Value = 3771.15

$ws.Range("A2007").Value = 'NELCO'
$ws.Range("F2007").Value = 881.2

$ws.Range("A2008").Value = 'NESCO'
$ws.Range("F2008").Value = 943.5

$ws.Range("A2009").Value = 'NEULANDLAB'
$ws.Range("F2009").Value = 9051.200000000001

$ws.Range("A2010").Value = 'NEXTMEDIA'
$ws.Range("F2010").Value = 7.33

$ws.Range("A2011").Value = 'NIBL'
$ws.Range("F2011").Value = 46.16

$ws.Range("A2012").Value = 'NRAIL'
$ws.Range("F2012").Value = 474.6

$ws.Range("A2013").Value = 'PARADEEP'
$ws.Range("F2013").Value = 94.48999999999999

$ws.Range("A2014").Value = 'PCBL'
$ws.Range("F2014").Value = 331.95

$ws.Range("A2015").Value = 'RAIN'
$ws.Range("F2015").Value = 172.76

$ws.Range("A2016").Value = 'RBA'
$ws.Range("F2016").Value = 110.49

$ws.Range("A2017").Value = 'RITCO'
$ws.Range("F2017").Value = 358.05

$ws.Range("A2018").Value = 'RPSGVENT'
$ws.Range("F2018").Value = 801.55

$ws.Range("A2019").Value = 'RRKABEL'
$ws.Range("F2019").Value = 1800.45

$ws.Range("A2020").Value = 'RUSTOMJEE'
$ws.Range("F2020").Value = 712.1

$ws.Range("A2021").Value = 'SAGCEM'
$ws.Range("F2021").Value = 248.6

$ws.Range("A2022").Value = '30/07/2024'

$ws.Range("A2023").Value = 'Buying Opportunity'
$ws.Range("B2023").Value = 'support Zone'
$ws.Range("C2023").Value = 'long buildup'
$ws.Range("D2023").Value = 'Short buildup'
$ws.Range("E2023").Value = 'FII ENTERING'

$ws.Range("A2024").Value = '5PAISA'
$ws.Range("B2024").Value = 'ALKALI'
$ws.Range("C2024").Value = 'BALRAMCHIN'
$ws.Range("E2024").Value = 'ASIANPAINT'
$ws.Range("F2024").Value = 491.45
$ws.Range("G2024").Value = 129.06
$ws.Range("H2024").Value = 480.85
$ws.Range("J2024").Value = 3084.45

$ws.Range("A2025").Value = 'AARVEEDEN'
$ws.Range("B2025").Value = 'ALLSEC'
$ws.Range("C2025").Value = 'HDFCLIFE'
$ws.Range("E2025").Value = 'DIXON'
$ws.Range("F2025").Value = 40.46
$ws.Range("G2025").Value = 1082.7
$ws.Range("H2025").Value = 715.5
$ws.Range("J2025").Value = 12106.45

$ws.Range("A2026").Value = 'ADANIENT'
$ws.Range("B2026").Value = 'CANFINHOME'
$ws.Range("C2026").Value = 'HINDCOPPER'
$ws.Range("E2026").Value = 'GODREJPROP'
$ws.Range("F2026").Value = 3169.4
$ws.Range("G2026").Value = 848.65
$ws.Range("H2026").Value = 321.6
$ws.Range("J2026").Value = 3219.55

$ws.Range("A2027").Value = 'ADANIPOWER'
$ws.Range("B2027").Value = 'DHANI'
$ws.Range("C2027").Value = 'JSWSTEEL'
$ws.Range("E2027").Value = 'GRANULES'
$ws.Range("F2027").Value = 734.45
$ws.Range("G2027").Value = 55.34
$ws.Range("H2027").Value = 928.25
$ws.Range("J2027").Value = 630.1

$ws.Range("A2028").Value = 'ADL'
$ws.Range("B2028").Value = 'INDIACEM'
$ws.Range("C2028").Value = 'MARUTI'
$ws.Range("E2028").Value = 'MARUTI'
$ws.Range("F2028").Value = 98.23
$ws.Range("G2028").Value = 364.25
$ws.Range("H2028").Value = 13115.8
$ws.Range("J2028").Value = 13115.8

$ws.Range("A2029").Value = 'AJANTPHARM'
$ws.Range("B2029").Value = 'INTELLECT'
$ws.Range("E2029").Value = 'POLYCAB'
$ws.Range("F2029").Value = 2695.1
$ws.Range("G2029").Value = 969.55
$ws.Range("J2029").Value = 6858.2

$ws.Range("A2030").Value = 'ANGELONE'
$ws.Range("B2030").Value = 'LANDMARK'
$ws.Range("F2030").Value = 2166.15
$ws.Range("G2030").Value = 628.55

$ws.Range("A2031").Value = 'ARROWGREEN'
$ws.Range("B2031").Value = 'MAHLOG'
$ws.Range("F2031").Value = 743.1
$ws.Range("G2031").Value = 514.2

$ws.Range("A2032").Value = 'ASIANPAINT'
$ws.Range("B2032").Value = 'MAITHANALL'
$ws.Range("F2032").Value = 3084.45
$ws.Range("G2032").Value = 1087.25

$ws.Range("A2033").Value = 'AUTOIND'
$ws.Range("B2033").Value = 'MALLCOM'
$ws.Range("F2033").Value = 144.58
$ws.Range("G2033").Value = 1389.25

$ws.Range("A2034").Value = 'AYMSYNTEX'
$ws.Range("B2034").Value = 'MHRIL'
$ws.Range("F2034").Value = 137
$ws.Range("G2034").Value = 456.85

$ws.Range("A2035").Value = 'BALPHARMA'
$ws.Range("B2035").Value = 'MTNL'
$ws.Range("F2035").Value = 121.22
$ws.Range("G2035").Value = 83.22

$ws.Range("A2036").Value = 'BSL'
$ws.Range("B2036").Value = 'ORICONENT'
$ws.Range("F2036").Value = 235.73
$ws.Range("G2036").Value = 41.62

$ws.Range("A2037").Value = 'CHAMBLFERT'
$ws.Range("B2037").Value = 'PGHH'
$ws.Range("F2037").Value = 518
$ws.Range("G2037").Value = 16904.55

$ws.Range("A2038").Value = 'CRAFTSMAN'
$ws.Range("B2038").Value = 'PLAZACABLE'
$ws.Range("F2038").Value = 5589.5
$ws.Range("G2038").Value = 86.06999999999999

$ws.Range("A2039").Value = 'DECCANCE'
$ws.Range("B2039").Value = 'RATNAVEER'
$ws.Range("F2039").Value = 719.8
$ws.Range("G2039").Value = 192.68

$ws.Range("A2040").Value = 'DIXON'
$ws.Range("F2040").Value = 12106.45

$ws.Range("A2041").Value = 'EDELWEISS'
$ws.Range("F2041").Value = 69.78

$ws.Range("A2042").Value = 'FIBERWEB'
$ws.Range("F2042").Value = 64.26000000000001

$ws.Range("A2043").Value = 'GAIL'
$ws.Range("F2043").Value = 240.97

$ws.Range("A2044").Value = 'GRANULES'
$ws.Range("F2044").Value = 630.1

$ws.Range("A2045").Value = 'GSFC'
$ws.Range("F2045").Value = 246.12

$ws.Range("A2046").Value = 'HMVL'
$ws.Range("F2046").Value = 105.69

$ws.Range("A2047").Value = 'HUBTOWN'
$ws.Range("F2047").Value = 247.27

$ws.Range("A2048").Value = 'INSPIRISYS'
$ws.Range("F2048").Value = 183.72

$ws.Range("A2049").Value = 'KANORICHEM'
$ws.Range("F2049").Value = 132.48

$ws.Range("A2050").Value = 'KERNEX'
$ws.Range("F2050").Value = 588.65

$ws.Range("A2051").Value = 'MARUTI'
$ws.Range("F2051").Value = 13115.8

$ws.Range("A2052").Value = 'MCL'
$ws.Range("F2052").Value = 38.36

$ws.Range("A2053").Value = 'POLYCAB'
$ws.Range("F2053").Value = 6858.2

$ws.Range("A2054").Value = 'PREMIERPOL'
$ws.Range("F2054").Value = 235.97

$ws.Range("A2055").Value = 'RAIN'
$ws.Range("F2055").Value = 173.97

$ws.Range("A2056").Value = '31/07/2024'
